$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows appended after the existing last row (131) for date 2025-11-05 (serial 45966).
$row132 = @(502.87600000000009, 1052.605, 275.60999999999996, 568.29200000000003, 277.39500000000004, 483.81099999999998, 480.23699999999997, 175.91299999999998, 25.12, 120.462, 186.88200000000001, 199.10999999999999, 905.01199999999983, 887.77300000000002, 843.45, 298.089, 234.46599999999998, 137.86600000000001, 23.07, 81.364000000000004, 62.75, 42.3, 86.18, 3.14)
$row133 = @(475.58600000000001, 401.46799999999996, 153.67500000000001, 59.007999999999996, 22.094000000000001, 215.505, 90.561999999999998, 257.38300000000004, 216.953, 148.61499999999998, 271.75300000000004, 486.733, 659.13300000000004, 453.49499999999995, 249.077, 532.44099999999992, 504.37400000000002, 55.431000000000004, 0, 95.512, 74.247, 40.289000000000001, 25.367000000000001, 69.843999999999994)

$ws.Cells.Item(132, 1).Value = 45966
$ws.Cells.Item(132, 2).Value = "四方坪站充电量(kw)"
for ($i = 0; $i -lt $row132.Length; $i++) {
    $ws.Cells.Item(132, 3 + $i).Value = $row132[$i]
}

$ws.Cells.Item(133, 1).Value = 45966
$ws.Cells.Item(133, 2).Value = "高岭站充电量(kw)"
for ($i = 0; $i -lt $row133.Length; $i++) {
    $ws.Cells.Item(133, 3 + $i).Value = $row133[$i]
}

# Match the column A date style and column C:Z number style used by prior rows.
$ws.Range("A132:A133").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("C132:Z133").NumberFormat = "0.00_);[Red]\(0.00\)"

[void]$ws.Range("H136").Select()
